# Weekly Fruta/Hortaliza update: a new price-report row for "Ajo" (Chino /
# Primera) is inserted before the existing row 257, pushing every
# subsequent record down by one row (old row 257 -> 258, ..., old row
# 304 -> 305). The sheet's used range grows from A1:R304 to A1:R305.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 257; Excel shifts rows 257:304 down
# to 258:305 and extends the used range accordingly.
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row with the new report entry.
$ws.Cells.Item(257, 1).Value = 8
$ws.Cells.Item(257, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(257, 3).Value = "Coquimbo"
$ws.Cells.Item(257, 4).Value = 44785
$ws.Cells.Item(257, 5).Value = 4
$ws.Cells.Item(257, 6).Value = 100112003
$ws.Cells.Item(257, 7).Value = "Ajo"
$ws.Cells.Item(257, 8).Value = "Chino"
$ws.Cells.Item(257, 9).Value = "Primera"
$ws.Cells.Item(257, 10).Value = 480
$ws.Cells.Item(257, 11).Value = 26000
$ws.Cells.Item(257, 12).Value = 26500
$ws.Cells.Item(257, 13).Value = 26250
$ws.Cells.Item(257, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(257, 15).Value = "China"
$ws.Cells.Item(257, 16).Value = 2625
$ws.Cells.Item(257, 17).Value = 10
$ws.Cells.Item(257, 18).Value = "Hortaliza"
